# update F.py usage_example.py 推薦功能要從980作為inital value
# Append two new log rows (9 and 10) that mirror rows 7 ("a1") and 8 ("b2"),
# extending the sheet's used range from A1:K8 to A1:K10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 - duplicate of row 7 ("a1")
$ws.Range("A9").Value = "a1"
$ws.Range("B9").Value = 0.8646729588508606
$ws.Range("C9").Value = 0.4265280067920685
$ws.Range("D9").Value = 0.8276968598365784
$ws.Range("E9").Value = 0.8786906003952026
$ws.Range("F9").Value = 0.8721588850021362
$ws.Range("G9").Value = 97.02021026611328
$ws.Range("H9").Value = 12.51860427856445
$ws.Range("I9").Value = 12.0368595123291
$ws.Range("J9").Value = 178.1204681396484
$ws.Range("K9").Value = 194.4346466064453

# Row 10 - duplicate of row 8 ("b2")
$ws.Range("A10").Value = "b2"
$ws.Range("B10").Value = 0.8646729588508606
$ws.Range("C10").Value = 0.4265280067920685
$ws.Range("D10").Value = 0.8276968598365784
$ws.Range("E10").Value = 0.8786906003952026
$ws.Range("F10").Value = 0.8721616268157959
$ws.Range("G10").Value = 97.02021026611328
$ws.Range("H10").Value = 12.52402591705322
$ws.Range("I10").Value = 12.03680610656738
$ws.Range("J10").Value = 178.1204681396484
$ws.Range("K10").Value = 194.4124450683594

# New row labels (A9, A10) use the same style as the existing label column
$ws.Range("A7:A8").Copy()
$ws.Range("A9:A10").PasteSpecial(-4122)
$excel.CutCopyMode = 0
